$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain text so values like
# "43.161.48" or "303.83" are not auto-converted to numbers/dates,
# matching the original inline-string cell type.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '43.161.48'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '2.371.93'
$ws.Range('E3').Value = '  +1.65%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '303.83'
$ws.Range('D6').Value = '95.76'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('D7').Value = '0.503'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '0.481'
$ws.Range('E9').Value = '  -3.24%  '
$ws.Range('D10').Value = '34.36'
$ws.Range('E10').Value = '  +0.32%  '
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').Value = '18.53'
$ws.Range('E13').Value = '  -3.73%  '
$ws.Range('E14').Value = '  +0.02%  '
$ws.Range('D15').Value = '2.737.13'
$ws.Range('E15').Value = '  +1.56%  '
$ws.Range('D16').Value = '2.373.32'
$ws.Range('E16').Value = '  +1.85%  '
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('D18').Value = '43.148.62'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('D19').Value = '12.01'
$ws.Range('E19').Value = '  -1.21%  '
$ws.Range('D20').Value = '6.31'
$ws.Range('E20').Value = '  +2.04%  '
$ws.Range('D21').Value = '0.0₃0886'
$ws.Range('E21').Value = '  -0.67%  '
$ws.Range('D22').Value = '68.17'
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').Value = '235.38'
$ws.Range('E23').Value = '  -0.77%  '
$ws.Range('E24').Value = '  -2.41%  '
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').Value = '24.60'
$ws.Range('E27').Value = '  -0.45%  '
$ws.Range('E28').Value = '  +15.08%  '
$ws.Range('E29').Value = '  +2.55%  '
$ws.Range('D30').Value = '32.36'
$ws.Range('E30').Value = '  +2.24%  '
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').Value = '5.02'
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('D33').Value = '17.64'
$ws.Range('E33').Value = '  -0.50%  '
$ws.Range('E34').Value = '  +1.98%  '
$ws.Range('E35').Value = '  +5.29%  '
$ws.Range('E36').Value = '  +1.92%  '
$ws.Range('D37').Value = '4.36'
$ws.Range('E37').Value = '  -0.84%  '
$ws.Range('D38').Value = '2.85'
$ws.Range('E38').Value = '  +3.74%  '
$ws.Range('D39').Value = '123.25'
$ws.Range('E39').Value = '  -11.05%  '
$ws.Range('D40').Value = '2.27'
$ws.Range('E40').Value = '  -1.51%  '
$ws.Range('E41').Value = '  -0.96%  '
$ws.Range('D42').Value = '21.09'
$ws.Range('E42').Value = '  -5.53%  '
$ws.Range('D43').Value = '1.938.63'
$ws.Range('E43').Value = '  +0.19%  '
$ws.Range('E44').Value = '  -0.26%  '
$ws.Range('D45').Value = '2.15'
$ws.Range('E45').Value = '  +4.35%  '
$ws.Range('E46').Value = '  -7.19%  '
$ws.Range('E47').Value = '  -1.00%  '
$ws.Range('D48').Value = '2.602.08'
$ws.Range('E48').Value = '  +1.58%  '
$ws.Range('E49').Value = '  +2.32%  '
$ws.Range('D50').Value = '71.91'
$ws.Range('E50').Value = '  -1.56%  '
$ws.Range('E51').Value = '  +0.96%  '

# Restore default (General/Normal) style so no new style index is
# introduced, matching the target workbook exactly.
$priceRange.Style = "Normal"
